# Modifications pour utiliser XGBClassifier et ajuster les prédictions
#
# This converts the workbook's previously regression-style price columns
# (C/D/E on the 'Valeurs réelles' sheet, and B/C/D on 'Prédictions') into
# classifier-style integer class labels, and renames the affected headers
# on 'Valeurs réelles' to reflect the new '_class' targets.

$wb = $excel.ActiveWorkbook

# --- Sheet 'Valeurs réelles' (sheet1): header renames + classification values ---
$ws1 = $wb.Worksheets.Item("Valeurs réelles")
$ws1.Range("C1").Value = "PRIX EXP POMME GOLDEN FRANCE 201/270G CAT.I PLATEAU 1RG_S+1_class"
$ws1.Range("D1").Value = "PRIX EXP POMME GOLDEN FRANCE 201/270G CAT.I PLATEAU 1RG_S+2_class"
$ws1.Range("E1").Value = "PRIX EXP POMME GOLDEN FRANCE 201/270G CAT.I PLATEAU 1RG_S+3_class"

$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 2
$ws1.Range("E2").Value = 2
$ws1.Range("C3").Value = 2
$ws1.Range("D3").Value = 2
$ws1.Range("E3").Value = 2
$ws1.Range("C4").Value = 2
$ws1.Range("D4").Value = 2
$ws1.Range("E4").Value = 2
$ws1.Range("C5").Value = 2
$ws1.Range("D5").Value = 2
$ws1.Range("E5").Value = 2
$ws1.Range("C6").Value = 2
$ws1.Range("D6").Value = 2
$ws1.Range("E6").Value = 2
$ws1.Range("C7").Value = 2
$ws1.Range("D7").Value = 2
$ws1.Range("E7").Value = 2
$ws1.Range("C8").Value = 2
$ws1.Range("D8").Value = 2
$ws1.Range("E8").Value = 2
$ws1.Range("C9").Value = 2
$ws1.Range("D9").Value = 2
$ws1.Range("E9").Value = 4
$ws1.Range("C10").Value = 2
$ws1.Range("D10").Value = 4
$ws1.Range("E10").Value = 0
$ws1.Range("C11").Value = 4
$ws1.Range("D11").Value = 0
$ws1.Range("E11").Value = 3
$ws1.Range("C12").Value = 0
$ws1.Range("D12").Value = 3
$ws1.Range("E12").Value = 4
$ws1.Range("C13").Value = 3
$ws1.Range("D13").Value = 4
$ws1.Range("E13").Value = 2
$ws1.Range("C14").Value = 4
$ws1.Range("D14").Value = 2
$ws1.Range("E14").Value = 0
$ws1.Range("C15").Value = 2
$ws1.Range("D15").Value = 0
$ws1.Range("E15").Value = 2
$ws1.Range("C16").Value = 0
$ws1.Range("D16").Value = 2
$ws1.Range("E16").Value = 1
$ws1.Range("C17").Value = 2
$ws1.Range("D17").Value = 1
$ws1.Range("E17").Value = 4
$ws1.Range("C18").Value = 1
$ws1.Range("D18").Value = 4
$ws1.Range("E18").Value = 0
$ws1.Range("C19").Value = 4
$ws1.Range("D19").Value = 0
$ws1.Range("E19").Value = 0
$ws1.Range("C20").Value = 0
$ws1.Range("D20").Value = 0
$ws1.Range("E20").Value = 1
$ws1.Range("C21").Value = 0
$ws1.Range("D21").Value = 1
$ws1.Range("E21").Value = 2
$ws1.Range("C22").Value = 1
$ws1.Range("D22").Value = 2
$ws1.Range("E22").Value = 4
$ws1.Range("C23").Value = 2
$ws1.Range("D23").Value = 4
$ws1.Range("E23").Value = 2
$ws1.Range("C24").Value = 4
$ws1.Range("D24").Value = 2
$ws1.Range("E24").Value = 2
$ws1.Range("C25").Value = 2
$ws1.Range("D25").Value = 2
$ws1.Range("E25").Value = 1
$ws1.Range("C26").Value = 2
$ws1.Range("D26").Value = 1
$ws1.Range("E26").Value = 2
$ws1.Range("C27").Value = 1
$ws1.Range("D27").Value = 2
$ws1.Range("E27").Value = 2
$ws1.Range("C28").Value = 2
$ws1.Range("D28").Value = 2
$ws1.Range("E28").Value = 2

# --- Sheet 'Prédictions' (sheet2): replace regression predictions with classifier outputs ---
$ws2 = $wb.Worksheets.Item("Prédictions")
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0
$ws2.Range("B3").Value = 0
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 0
$ws2.Range("B4").Value = 0
$ws2.Range("C4").Value = 0
$ws2.Range("D4").Value = 0
$ws2.Range("B5").Value = 0
$ws2.Range("C5").Value = 0
$ws2.Range("D5").Value = 0
$ws2.Range("B6").Value = 0
$ws2.Range("C6").Value = 0
$ws2.Range("D6").Value = 0
$ws2.Range("B7").Value = 0
$ws2.Range("C7").Value = 0
$ws2.Range("D7").Value = 0
$ws2.Range("B8").Value = 2
$ws2.Range("C8").Value = 0
$ws2.Range("D8").Value = 0
$ws2.Range("B9").Value = 0
$ws2.Range("C9").Value = 0
$ws2.Range("D9").Value = 0
$ws2.Range("B10").Value = 2
$ws2.Range("C10").Value = 0
$ws2.Range("D10").Value = 0
$ws2.Range("B11").Value = -2
$ws2.Range("C11").Value = 0
$ws2.Range("D11").Value = 2
$ws2.Range("B12").Value = 0
$ws2.Range("C12").Value = 0
$ws2.Range("D12").Value = 0
$ws2.Range("B13").Value = 0
$ws2.Range("C13").Value = 0
$ws2.Range("D13").Value = 0
$ws2.Range("B14").Value = 0
$ws2.Range("C14").Value = 0
$ws2.Range("D14").Value = 0
$ws2.Range("B15").Value = 0
$ws2.Range("C15").Value = 0
$ws2.Range("D15").Value = 0
$ws2.Range("B16").Value = 0
$ws2.Range("C16").Value = 0
$ws2.Range("D16").Value = 0
$ws2.Range("B17").Value = 0
$ws2.Range("C17").Value = 0
$ws2.Range("D17").Value = -2
$ws2.Range("B18").Value = 0
$ws2.Range("C18").Value = 0
$ws2.Range("D18").Value = 0
$ws2.Range("B19").Value = 0
$ws2.Range("C19").Value = 0
$ws2.Range("D19").Value = 0
$ws2.Range("B20").Value = 0
$ws2.Range("C20").Value = 0
$ws2.Range("D20").Value = 0
$ws2.Range("B21").Value = 0
$ws2.Range("C21").Value = 0
$ws2.Range("D21").Value = 0
$ws2.Range("B22").Value = 0
$ws2.Range("C22").Value = 0
$ws2.Range("D22").Value = 0
$ws2.Range("B23").Value = 0
$ws2.Range("C23").Value = 0
$ws2.Range("D23").Value = 1
$ws2.Range("B24").Value = 0
$ws2.Range("C24").Value = 0
$ws2.Range("D24").Value = 0
$ws2.Range("B25").Value = 0
$ws2.Range("C25").Value = 0
$ws2.Range("D25").Value = -2
$ws2.Range("B26").Value = 0
$ws2.Range("C26").Value = 0
$ws2.Range("D26").Value = 2
$ws2.Range("B27").Value = 0
$ws2.Range("C27").Value = 0
$ws2.Range("D27").Value = 0
$ws2.Range("B28").Value = 0
$ws2.Range("C28").Value = 0
$ws2.Range("D28").Value = -1
